$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2415.5938
$ws.Range("I17").Value = 281.5
$ws.Range("J17").Value = 2484.4355
$ws.Range("K17").Value = 844.5
$ws.Range("L17").Value = 7453.306500000001
$ws.Range("M17").Value = -676.5
$ws.Range("N17").Value = -7789.306500000001
$ws.Range("H32").Value = 23810462
$ws.Range("I32").Value = 66666944
$ws.Range("J32").Value = 1305.5555
$ws.Range("K32").Value = 66666944
$ws.Range("L32").Value = 1305.5555
$ws.Range("M32").Value = -66666618
$ws.Range("N32").Value = -1957.5555
$ws.Range("H106").Value = 23405.695
$ws.Range("I106").Value = 1576.1666
$ws.Range("J106").Value = 67064.75
$ws.Range("K106").Value = 1576.1666
$ws.Range("L106").Value = 67064.75
$ws.Range("M106").Value = -945.1666
$ws.Range("N106").Value = -68326.75
$ws.Range("H108").Value = 45207.332
$ws.Range("J108").Value = 45207.332
$ws.Range("L108").Value = 45207.332
$ws.Range("N108").Value = -52887.332
$ws.Range("H117").Value = 48718.4
$ws.Range("J117").Value = 48718.4
$ws.Range("L117").Value = 48718.4
$ws.Range("N117").Value = -57896.4
$ws.Range("H120").Value = 49702
$ws.Range("J120").Value = 49702
$ws.Range("L120").Value = 49702
$ws.Range("N120").Value = -59378
$ws.Range("H121").Value = 1311
$ws.Range("I121").Value = 298.33334
$ws.Range("J121").Value = 2830
$ws.Range("K121").Value = 895.0000200000001
$ws.Range("L121").Value = 8490
$ws.Range("M121").Value = 851.9999799999999
$ws.Range("N121").Value = -11984
$ws.Range("H124").Value = 43698
$ws.Range("J124").Value = 43698
$ws.Range("L124").Value = 43698
$ws.Range("N124").Value = -53518
$ws.Range("H126").Value = 46772
$ws.Range("J126").Value = 46772
$ws.Range("L126").Value = 46772
$ws.Range("N126").Value = -56652
$ws.Range("H128").Value = 49580.75
$ws.Range("J128").Value = 49580.75
$ws.Range("L128").Value = 49580.75
$ws.Range("N128").Value = -59540.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 36856.75
$ws.Range("J109").Value = 36856.75
$ws.Range("L109").Value = 36856.75
$ws.Range("N109").Value = -39630.75
$ws.Range("H122").Value = 1778.7407
$ws.Range("I122").Value = 1882.35
$ws.Range("J122").Value = 1482.7142
$ws.Range("K122").Value = 5647.049999999999
$ws.Range("L122").Value = 4448.142599999999
$ws.Range("M122").Value = -3197.049999999999
$ws.Range("N122").Value = -9348.142599999999
$ws.Range("H123").Value = 49710
$ws.Range("J123").Value = 49710
$ws.Range("L123").Value = 49710
$ws.Range("N123").Value = -59510
$ws.Range("H125").Value = 35916.855
$ws.Range("J125").Value = 35916.855
$ws.Range("L125").Value = 35916.855
$ws.Range("N125").Value = -45756.855
$ws.Range("H131").Value = 46401.75
$ws.Range("J131").Value = 46401.75
$ws.Range("L131").Value = 46401.75
$ws.Range("N131").Value = -56481.75
$ws.Range("H139").Value = 49433.168
$ws.Range("J139").Value = 49433.168
$ws.Range("L139").Value = 49433.168
$ws.Range("N139").Value = -59713.168

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 3500
$ws.Range("J38").Value = 3500
$ws.Range("L38").Value = 3500
$ws.Range("N38").Value = -4332
$ws.Range("H99").Value = 2333.8696
$ws.Range("I99").Value = 2159.95
$ws.Range("J99").Value = 3493.3333
$ws.Range("K99").Value = 2159.95
$ws.Range("L99").Value = 3493.3333
$ws.Range("M99").Value = -661.9499999999998
$ws.Range("N99").Value = -6489.3333
$ws.Range("H102").Value = 15751.5
$ws.Range("I102").Value = 8670.666999999999
$ws.Range("K102").Value = 8670.666999999999
$ws.Range("M102").Value = -5425.666999999999
$ws.Range("H117").Value = 45996
$ws.Range("J117").Value = 45996
$ws.Range("L117").Value = 45996
$ws.Range("N117").Value = -55174
$ws.Range("H125").Value = 50780
$ws.Range("J125").Value = 50780
$ws.Range("L125").Value = 50780
$ws.Range("N125").Value = -60620
$ws.Range("H130").Value = 49889.5
$ws.Range("J130").Value = 49889.5
$ws.Range("L130").Value = 49889.5
$ws.Range("N130").Value = -59929.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49883.5
$ws.Range("J20").Value = 49883.5
$ws.Range("L20").Value = 49883.5
$ws.Range("N20").Value = -50355.5
$ws.Range("H30").Value = 49883.5
$ws.Range("J30").Value = 49883.5
$ws.Range("L30").Value = 49883.5
$ws.Range("N30").Value = -50065.5
$ws.Range("H100").Value = 44617.332
$ws.Range("J100").Value = 44617.332
$ws.Range("L100").Value = 44617.332
$ws.Range("N100").Value = -46781.332
$ws.Range("H128").Value = 49883.5
$ws.Range("J128").Value = 49883.5
$ws.Range("L128").Value = 49883.5
$ws.Range("N128").Value = -59843.5
$ws.Range("H141").Value = 10000
$ws.Range("J141").Value = 10000
$ws.Range("L141").Value = 10000
$ws.Range("N141").Value = -20360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2470.2
$ws.Range("I102").Value = 2522.4443
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2522.4443
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -900.4443000000001
$ws.Range("N102").Value = -5244
$ws.Range("H110").Value = 51999
$ws.Range("J110").Value = 51999
$ws.Range("L110").Value = 51999
$ws.Range("N110").Value = -60179
$ws.Range("H124").Value = 38415
$ws.Range("I124").Value = 31709
$ws.Range("J124").Value = 41768
$ws.Range("K124").Value = 31709
$ws.Range("L124").Value = 41768
$ws.Range("M124").Value = -26799
$ws.Range("N124").Value = -51588
$ws.Range("H130").Value = 45924
$ws.Range("J130").Value = 45924
$ws.Range("L130").Value = 45924
$ws.Range("N130").Value = -55964
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 3000
$ws.Range("J38").Value = 3000
$ws.Range("L38").Value = 3000
$ws.Range("N38").Value = -3820
$ws.Range("H47").Value = 26136.428
$ws.Range("J47").Value = 26136.428
$ws.Range("L47").Value = 26136.428
$ws.Range("N47").Value = -27116.428
$ws.Range("H52").Value = 26136.428
$ws.Range("J52").Value = 26136.428
$ws.Range("L52").Value = 26136.428
$ws.Range("N52").Value = -26602.428
$ws.Range("H111").Value = 46249
$ws.Range("J111").Value = 46249
$ws.Range("L111").Value = 46249
$ws.Range("N111").Value = -54429
$ws.Range("H121").Value = 43420
$ws.Range("J121").Value = 43420
$ws.Range("L121").Value = 43420
$ws.Range("N121").Value = -46914
$ws.Range("H124").Value = 47693
$ws.Range("J124").Value = 47693
$ws.Range("L124").Value = 47693
$ws.Range("N124").Value = -57513
$ws.Range("H125").Value = 48711
$ws.Range("J125").Value = 48711
$ws.Range("L125").Value = 48711
$ws.Range("N125").Value = -58551
$ws.Range("H127").Value = 50715
$ws.Range("J127").Value = 50715
$ws.Range("L127").Value = 50715
$ws.Range("N127").Value = -60635
$ws.Range("H128").Value = 48421
$ws.Range("J128").Value = 48421
$ws.Range("L128").Value = 48421
$ws.Range("N128").Value = -58381
$ws.Range("H130").Value = 48210.5
$ws.Range("J130").Value = 48210.5
$ws.Range("L130").Value = 48210.5
$ws.Range("N130").Value = -58250.5
$ws.Range("H134").Value = 44249.75
$ws.Range("J134").Value = 44249.75
$ws.Range("L134").Value = 44249.75
$ws.Range("N134").Value = -54389.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 40674
$ws.Range("J119").Value = 40674
$ws.Range("L119").Value = 40674
$ws.Range("N119").Value = -50350
$ws.Range("H131").Value = 50136
$ws.Range("J131").Value = 50136
$ws.Range("L131").Value = 50136
$ws.Range("N131").Value = -60216
